$wb = $excel.ActiveWorkbook

# --- Overall sheet, row 2 ---
$wsOverall = $wb.Worksheets.Item("Overall")
$wsOverall.Range("B2").Value = 76
$wsOverall.Range("C2").Value = 43
$wsOverall.Range("D2").Value = 1.1072865120067301
$wsOverall.Range("E2").Value = 0.55175438596491255
$wsOverall.Range("F2").Value = 1.3102694042143173
$wsOverall.Range("G2").Value = 42
$wsOverall.Range("H2").Value = 31
$wsOverall.Range("I2").Value = 73
$wsOverall.Range("J2").Value = 537
$wsOverall.Range("K2").Value = 23

# --- Zones sheet, rows 2-14 ---
$wsZones = $wb.Worksheets.Item("Zones")
# row 2
$wsZones.Range("B2").Value = 5
$wsZones.Range("C2").Value = 8
$wsZones.Range("D2").Value = 0.98666666666666691
$wsZones.Range("E2").Value = 0.23333333333333339
$wsZones.Range("F2").Value = 1.0404761904761908
# row 3
$wsZones.Range("B3").Value = 0
$wsZones.Range("C3").Value = 4
$wsZones.Range("D3").Value = 1.8404761904761904
$wsZones.Range("F3").Value = 1.8404761904761904
# row 4
$wsZones.Range("B4").Value = 8
$wsZones.Range("C4").Value = 5
$wsZones.Range("D4").Value = 1.6714285714285713
$wsZones.Range("E4").Value = 0.6333333333333333
$wsZones.Range("F4").Value = 1.9545454545454546
# row 5
$wsZones.Range("B5").Value = 5
$wsZones.Range("C5").Value = 1
$wsZones.Range("D5").Value = 0.82222222222222197
$wsZones.Range("E5").Value = 0.33333333333333348
$wsZones.Range("F5").Value = 0.91999999999999971
# row 6
$wsZones.Range("B6").Value = 4
$wsZones.Range("C6").Value = 0
$wsZones.Range("D6").Value = 0.50416666666666665
$wsZones.Range("E6").Value = 0.50833333333333375
$wsZones.Range("F6").Value = 0.49999999999999956
# row 7
$wsZones.Range("B7").Value = 13
$wsZones.Range("C7").Value = 2
$wsZones.Range("D7").Value = 0.67395833333333344
$wsZones.Range("E7").Value = 0.45714285714285724
$wsZones.Range("F7").Value = 0.84259259259259267
# row 8
$wsZones.Range("B8").Value = 8
$wsZones.Range("C8").Value = 2
$wsZones.Range("D8").Value = 0.63000000000000012
$wsZones.Range("E8").Value = 0.41333333333333344
$wsZones.Range("F8").Value = 0.84666666666666668
# row 9
$wsZones.Range("B9").Value = 2
$wsZones.Range("C9").Value = 2
$wsZones.Range("D9").Value = 0.43958333333333333
$wsZones.Range("E9").Value = 0.38333333333333341
$wsZones.Range("F9").Value = 0.45833333333333326
# row 10
$wsZones.Range("B10").Value = 1
$wsZones.Range("C10").Value = 6
$wsZones.Range("D10").Value = 1.6621212121212119
$wsZones.Range("F10").Value = 1.6621212121212119
# row 11
$wsZones.Range("C11").Value = 4
$wsZones.Range("D11").Value = 1.4051282051282052
$wsZones.Range("E11").Value = 1.6944444444444444
$wsZones.Range("F11").Value = 1.3183333333333336
# row 12
$wsZones.Range("B12").Value = 8
$wsZones.Range("C12").Value = 0
$wsZones.Range("D12").Value = 0.3708333333333334
$wsZones.Range("E12").Value = 0.38333333333333336
$wsZones.Range("F12").Value = 0.33333333333333348
# row 13
$wsZones.Range("B13").Value = 7
$wsZones.Range("C13").Value = 4
$wsZones.Range("D13").Value = 1.2013888888888891
$wsZones.Range("E13").Value = 0.61333333333333351
$wsZones.Range("F13").Value = 1.6214285714285717
# row 14
$wsZones.Range("B14").Value = 11
$wsZones.Range("C14").Value = 5
$wsZones.Range("D14").Value = 1.3870370370370371
$wsZones.Range("E14").Value = 0.33333333333333331
$wsZones.Range("F14").Value = 1.597777777777778

Write-Output "edits applied"
